$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.20011764684948
$ws.Range("C2").Value = 5.549510265386702
$ws.Range("D2").Value = 13.59115401781157
$ws.Range("E2").Value = 13.88680011581545
$ws.Range("G2").Value = 3.7188001584187
$ws.Range("I2").Value = 31.0040753684009
$ws.Range("J2").Value = 8.431817319517179
$ws.Range("K2").Value = 12.41729235547219
$ws.Range("L2").Value = 12.54101177934104
$ws.Range("M2").Value = 17.05710824708649
$ws.Range("N2").Value = 22.72436872851446
$ws.Range("O2").Value = 32.53314787328899
$ws.Range("B3").Value = 15.05679776130088
$ws.Range("C3").Value = 5.460530465355505
$ws.Range("D3").Value = 13.59330448255508
$ws.Range("E3").Value = 13.9114036168101
$ws.Range("G3").Value = 3.720720094090111
$ws.Range("I3").Value = 31.08948710429331
$ws.Range("J3").Value = 8.435105540830625
$ws.Range("K3").Value = 12.31728889375758
$ws.Range("L3").Value = 12.55427568381321
$ws.Range("M3").Value = 17.04602917981959
$ws.Range("N3").Value = 22.78537168138178
$ws.Range("O3").Value = 32.61261227130181
$ws.Range("B4").Value = 14.97094253569502
$ws.Range("C4").Value = 5.40432690750312
$ws.Range("D4").Value = 13.5968356271749
$ws.Range("E4").Value = 13.92792625868897
$ws.Range("G4").Value = 3.721962626316204
$ws.Range("I4").Value = 31.14633199124787
$ws.Range("J4").Value = 8.437246457370213
$ws.Range("K4").Value = 12.25763757101857
$ws.Range("L4").Value = 12.56387412749271
$ws.Range("M4").Value = 17.04149912648062
$ws.Range("N4").Value = 22.82460548094216
$ws.Range("O4").Value = 32.66649561385207
$ws.Range("B5").Value = 14.93653029067612
$ws.Range("C5").Value = 5.38104202536075
$ws.Range("D5").Value = 13.59883169017946
$ws.Range("E5").Value = 13.93501596980867
$ws.Range("G5").Value = 3.722485031669145
$ws.Range("I5").Value = 31.17060377681535
$ws.Range("J5").Value = 8.438149641115801
$ws.Range("K5").Value = 12.23379208598873
$ws.Range("L5").Value = 12.56815187669031
$ws.Range("M5").Value = 17.04022735863803
$ws.Range("N5").Value = 22.84104186825291
$ws.Range("O5").Value = 32.68973341852733
$ws.Range("B6").Value = 14.93085187307594
$ws.Range("C6").Value = 5.377152944352733
$ws.Range("D6").Value = 13.59919681258123
$ws.Range("E6").Value = 13.93621476435192
$ws.Range("G6").Value = 3.722572748326614
$ws.Range("I6").Value = 31.1747009555503
$ws.Range("J6").Value = 8.438301473275221
$ws.Range("K6").Value = 12.22986116383554
$ws.Range("L6").Value = 12.56888433428636
$ws.Range("M6").Value = 17.04005094452068
$ws.Range("N6").Value = 22.84379823111426
$ws.Range("O6").Value = 32.69366932582376
$ws.Range("B7").Value = 14.97047606819171
$ws.Range("C7").Value = 5.40401440618691
$ws.Range("D7").Value = 13.59686028971021
$ws.Range("E7").Value = 13.9280204284386
$ws.Range("G7").Value = 3.721969606551315
$ws.Range("I7").Value = 31.14665484647834
$ws.Range("J7").Value = 8.43725851343652
$ws.Range("K7").Value = 12.25731407893952
$ws.Range("L7").Value = 12.56393033479665
$ws.Range("M7").Value = 17.04147964624728
$ws.Range("N7").Value = 22.82482533100001
$ws.Range("O7").Value = 32.666803825434
$ws.Range("B8").Value = 15.15027754649878
$ws.Range("C8").Value = 5.519163561636234
$ws.Range("D8").Value = 13.5914374218917
$ws.Range("E8").Value = 13.89498983861921
$ws.Range("G8").Value = 3.719448964488588
$ws.Range("I8").Value = 31.03261186785434
$ws.Range("J8").Value = 8.432925839442277
$ws.Range("K8").Value = 12.3824614640812
$ws.Range("L8").Value = 12.54528369140397
$ws.Range("M8").Value = 17.05281827475467
$ws.Range("N8").Value = 22.74503434368988
$ws.Range("O8").Value = 32.55949010289918
$ws.Range("B9").Value = 15.51827872163138
$ws.Range("C9").Value = 5.731965249203537
$ws.Range("D9").Value = 13.59828569016984
$ws.Range("E9").Value = 13.84142840134471
$ws.Range("G9").Value = 3.715009038743737
$ws.Range("I9").Value = 30.84389154014965
$ws.Range("J9").Value = 8.425393306769394
$ws.Range("K9").Value = 12.64074509365187
$ws.Range("L9").Value = 12.52023031701288
$ws.Range("M9").Value = 17.09295318768033
$ws.Range("N9").Value = 22.60261115565569
$ws.Range("O9").Value = 32.38947053248998
$ws.Range("B10").Value = 15.79585588989691
$ws.Range("C10").Value = 5.879686373723009
$ws.Range("D10").Value = 13.6138857323613
$ws.Range("E10").Value = 13.808879267801
$ws.Range("G10").Value = 3.712050537971932
$ws.Range("I10").Value = 30.7265139328924
$ws.Range("J10").Value = 8.42044163329844
$ws.Range("K10").Value = 12.83693908008828
$ws.Range("L10").Value = 12.50880452407565
$ws.Range("M10").Value = 17.13316435406264
$ws.Range("N10").Value = 22.50645225603417
$ws.Range("O10").Value = 32.28922554491638
$ws.Range("B11").Value = 15.92322481009515
$ws.Range("C11").Value = 5.944880349557207
$ws.Range("D11").Value = 13.62325419728311
$ws.Range("E11").Value = 13.79554190542301
$ws.Range("G11").Value = 3.710769873314645
$ws.Range("I11").Value = 30.67773301989374
$ws.Range("J11").Value = 8.418314389595043
$ws.Range("K11").Value = 12.92728062937042
$ws.Range("L11").Value = 12.50511363892453
$ws.Range("M11").Value = 17.15374115617447
$ws.Range("N11").Value = 22.46453042794039
$ws.Range("O11").Value = 32.24898217819283
$ws.Range("B12").Value = 15.97157420888083
$ws.Range("C12").Value = 5.969268675728989
$ws.Range("D12").Value = 13.62712612754268
$ws.Range("E12").Value = 13.79070211117003
$ws.Range("G12").Value = 3.710294239788485
$ws.Range("I12").Value = 30.65992439739761
$ws.Range("J12").Value = 8.417526793291112
$ws.Range("K12").Value = 12.96162135907836
$ws.Range("L12").Value = 12.50393182899784
$ws.Range("M12").Value = 17.16185732671657
$ws.Range("N12").Value = 22.44891631250088
$ws.Range("O12").Value = 32.2345138383295
$ws.Range("B13").Value = 15.96115671110119
$ws.Range("C13").Value = 5.964029687784882
$ws.Range("D13").Value = 13.62627785556704
$ws.Range("E13").Value = 13.79173508159722
$ws.Range("G13").Value = 3.710396261902477
$ws.Range("I13").Value = 30.6637302861449
$ws.Range("J13").Value = 8.417695619304673
$ws.Range("K13").Value = 12.95422009616589
$ws.Range("L13").Value = 12.50417676665267
$ws.Range("M13").Value = 17.16009501760518
$ws.Range("N13").Value = 22.45226751172395
$ws.Range("O13").Value = 32.23759556076223
$ws.Range("B14").Value = 15.92720041809797
$ws.Range("C14").Value = 5.946892848724454
$ws.Range("D14").Value = 13.62356625849704
$ws.Range("E14").Value = 13.79513951141318
$ws.Range("G14").Value = 3.710730556001097
$ws.Range("I14").Value = 30.67625459299888
$ws.Range("J14").Value = 8.418249234349503
$ws.Range("K14").Value = 12.93010340238514
$ws.Range("L14").Value = 12.50501208968624
$ws.Range("M14").Value = 17.15440240697913
$ws.Range("N14").Value = 22.46324062548782
$ws.Range("O14").Value = 32.24777640627616
$ws.Range("B15").Value = 15.9064153258017
$ws.Range("C15").Value = 5.936356760622165
$ws.Range("D15").Value = 13.62194748477667
$ws.Range("E15").Value = 13.7972522564734
$ws.Range("G15").Value = 3.710936534022883
$ws.Range("I15").Value = 30.68401252083591
$ws.Range("J15").Value = 8.418590674330199
$ws.Range("K15").Value = 12.9153474139647
$ws.Range("L15").Value = 12.50555183367558
$ws.Range("M15").Value = 17.15095760403331
$ws.Range("N15").Value = 22.46999590110159
$ws.Range("O15").Value = 32.25411287854828
$ws.Range("B16").Value = 15.78755066979192
$ws.Range("C16").Value = 5.875384501411037
$ws.Range("D16").Value = 13.61331899763187
$ws.Range("E16").Value = 13.80978041945932
$ws.Range("G16").Value = 3.712135539969272
$ws.Range("I16").Value = 30.7297947459179
$ws.Range("J16").Value = 8.420583168734339
$ws.Range("K16").Value = 12.83105474066478
$ws.Range("L16").Value = 12.50907598289825
$ws.Range("M16").Value = 17.13186521762707
$ws.Range("N16").Value = 22.50922849876882
$ws.Range("O16").Value = 32.29196340761789
$ws.Range("B17").Value = 15.71488345902981
$ws.Range("C17").Value = 5.837458496178013
$ws.Range("D17").Value = 13.6086059546776
$ws.Range("E17").Value = 13.81784202188702
$ws.Range("G17").Value = 3.71288775109062
$ws.Range("I17").Value = 30.75906265043788
$ws.Range("J17").Value = 8.421837539098421
$ws.Range("K17").Value = 12.77960462484911
$ws.Range("L17").Value = 12.51162328377135
$ws.Range("M17").Value = 17.12073477692695
$ws.Range("N17").Value = 22.53376210371479
$ws.Range("O17").Value = 32.31655625138627
$ws.Range("B18").Value = 15.67319292239407
$ws.Range("C18").Value = 5.815456342506222
$ws.Range("D18").Value = 13.60610913951182
$ws.Range("E18").Value = 13.82261718131778
$ws.Range("G18").Value = 3.713326540030313
$ws.Range("I18").Value = 30.77633113466186
$ws.Range("J18").Value = 8.422570817656123
$ws.Range("K18").Value = 12.75011625879274
$ws.Range("L18").Value = 12.51323030018272
$ws.Range("M18").Value = 17.11454814819203
$ws.Range("N18").Value = 22.54804468441684
$ws.Range("O18").Value = 32.331205762024
$ws.Range("B19").Value = 15.65909661207325
$ws.Range("C19").Value = 5.807974823703026
$ws.Range("D19").Value = 13.60530058557486
$ws.Range("E19").Value = 13.82425774386133
$ws.Range("G19").Value = 3.713476161847161
$ws.Range("I19").Value = 30.78225254647249
$ws.Range("J19").Value = 8.422821121708557
$ws.Range("K19").Value = 12.74015074388486
$ws.Range("L19").Value = 12.5137987999591
$ws.Range("M19").Value = 17.11249056913593
$ws.Range("N19").Value = 22.5529100109142
$ws.Range("O19").Value = 32.33625244499964
$ws.Range("B20").Value = 15.72260835857545
$ws.Range("C20").Value = 5.841515334587174
$ws.Range("D20").Value = 13.60908553441562
$ws.Range("E20").Value = 13.816969536745
$ws.Range("G20").Value = 3.712807042105673
$ws.Range("I20").Value = 30.7559020769799
$ws.Range("J20").Value = 8.421702788724419
$ws.Range("K20").Value = 12.78507096767185
$ws.Range("L20").Value = 12.51133744137808
$ws.Range("M20").Value = 17.1218973785576
$ws.Range("N20").Value = 22.53113271822544
$ws.Range("O20").Value = 32.31388609893229
$ws.Range("B21").Value = 15.93717133384414
$ws.Range("C21").Value = 5.951934558091945
$ws.Range("D21").Value = 13.62435393924347
$ws.Range("E21").Value = 13.79413383168322
$ws.Range("G21").Value = 3.710632112964949
$ws.Range("I21").Value = 30.67255788883613
$ws.Range("J21").Value = 8.418086137753967
$ws.Range("K21").Value = 12.93718373504914
$ws.Range("L21").Value = 12.50476088381626
$ws.Range("M21").Value = 17.15606570220952
$ws.Range("N21").Value = 22.46001048772164
$ws.Range("O21").Value = 32.24476511915293
$ws.Range("B22").Value = 16.07806737740689
$ws.Range("C22").Value = 6.022351226453092
$ws.Range("D22").Value = 13.63622169872484
$ws.Range("E22").Value = 13.78043768909215
$ws.Range("G22").Value = 3.709265011608238
$ws.Range("I22").Value = 30.62195592120647
$ws.Range("J22").Value = 8.415827009489515
$ws.Range("K22").Value = 13.03734561814819
$ws.Range("L22").Value = 12.50172043369478
$ws.Range("M22").Value = 17.18028408148919
$ws.Range("N22").Value = 22.41504745314741
$ws.Range("O22").Value = 32.20408436205675
$ws.Range("B23").Value = 16.00282073120122
$ws.Range("C23").Value = 5.984931964099007
$ws.Range("D23").Value = 13.62971566203204
$ws.Range("E23").Value = 13.78763536295043
$ws.Range("G23").Value = 3.709989702031936
$ws.Range("I23").Value = 30.64860918205417
$ws.Range("J23").Value = 8.41702320521642
$ws.Range("K23").Value = 12.98382775935267
$ws.Range("L23").Value = 12.50322838363559
$ws.Range("M23").Value = 17.16718707169763
$ws.Range("N23").Value = 22.4389064128229
$ws.Range("O23").Value = 32.22538516002868
$ws.Range("B24").Value = 15.71911565800577
$ws.Range("C24").Value = 5.839681853823557
$ws.Range("D24").Value = 13.60886805347299
$ws.Range("E24").Value = 13.81736354979787
$ws.Range("G24").Value = 3.71284351090779
$ws.Range("I24").Value = 30.75732959550305
$ws.Range("J24").Value = 8.421763671598184
$ws.Range("K24").Value = 12.78259934861728
$ws.Range("L24").Value = 12.51146622666954
$ws.Range("M24").Value = 17.12137110439836
$ws.Range("N24").Value = 22.53232090920515
$ws.Range("O24").Value = 32.31509168389536
$ws.Range("B25").Value = 15.41729886372145
$ws.Range("C25").Value = 5.6758631818629
$ws.Range("D25").Value = 13.5945697838126
$ws.Range("E25").Value = 13.85472123578893
$ws.Range("G25").Value = 3.716156629601233
$ws.Range("I25").Value = 30.89120836745411
$ws.Range("J25").Value = 8.427328401456196
$ws.Range("K25").Value = 12.56963871574635
$ws.Range("L25").Value = 12.5257790968221
$ws.Range("M25").Value = 17.0801977791175
$ws.Range("N25").Value = 22.63964531285466
$ws.Range("O25").Value = 32.43113486040937
